$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force column D to text so numeric-looking price strings (e.g. "1.001")
# are stored as literal text instead of being parsed as numbers.
$priceRange = $ws.Range("D2:D51")
$priceRange.NumberFormat = "@"

$ws.Range("D2").Value = "29.146.05"
$ws.Range("D3").Value = "1.825.41"
$ws.Range("D4").Value = "0.9993"
$ws.Range("D5").Value = "241.86"
$ws.Range("D6").Value = "0.6175"
$ws.Range("D7").Value = "1.001"
$ws.Range("D8").Value = "0.07351"
$ws.Range("D9").Value = "0.2889"
$ws.Range("D10").Value = "22.95"
$ws.Range("D11").Value = "0.07677"
$ws.Range("D12").Value = "1.808.19"
$ws.Range("D13").Value = "4.961"
$ws.Range("D14").Value = "0.6623"
$ws.Range("D15").Value = "81.77"
$ws.Range("D16").Value = "0.000008955"
$ws.Range("D17").Value = "5.859"
$ws.Range("D18").Value = "29.119.77"
$ws.Range("D19").Value = "2.054.30"
$ws.Range("D20").Value = "237.89"
$ws.Range("D21").Value = "12.49"
$ws.Range("D22").Value = "1.001"
$ws.Range("D23").Value = "7.143"
$ws.Range("D25").Value = "158.38"
$ws.Range("D26").Value = "0.1410"
$ws.Range("D27").Value = "8.452"
$ws.Range("D28").Value = "17.66"
$ws.Range("D30").Value = "0.05591"
$ws.Range("D31").Value = "4.098"
$ws.Range("D32").Value = "4.108"
$ws.Range("D33").Value = "1.206"
$ws.Range("D34").Value = "1.825"
$ws.Range("D35").Value = "0.7350"
$ws.Range("D36").Value = "1.132"
$ws.Range("D37").Value = "2.623"
$ws.Range("D38").Value = "2.834"
$ws.Range("D39").Value = "1.208.47"
$ws.Range("D40").Value = "0.01762"
$ws.Range("D41").Value = "6.389"
$ws.Range("D42").Value = "0.8932"
$ws.Range("D43").Value = "1.001"
$ws.Range("D44").Value = "100.71"
$ws.Range("D45").Value = "1.957.86"
$ws.Range("D46").Value = "64.77"
$ws.Range("D47").Value = "0.00000000122"
$ws.Range("D48").Value = "0.5076"
$ws.Range("D49").Value = "9.070"
$ws.Range("D50").Value = "0.3998"
$ws.Range("D51").Value = "0.05795"
$ws.Range("E2").Value = "  +0.20%  "
$ws.Range("E3").Value = "  -0.48%  "
$ws.Range("E4").Value = "  +0.00%  "
$ws.Range("E5").Value = "  -0.77%  "
$ws.Range("E6").Value = "  -1.52%  "
$ws.Range("E7").Value = "  -0.05%  "
$ws.Range("E8").Value = "  -2.18%  "
$ws.Range("E9").Value = "  -1.16%  "
$ws.Range("E10").Value = "  -0.94%  "
$ws.Range("E11").Value = "  -0.08%  "
$ws.Range("E12").Value = "  -1.00%  "
$ws.Range("E13").Value = "  -0.81%  "
$ws.Range("E14").Value = "  -0.80%  "
$ws.Range("E15").Value = "  -1.09%  "
$ws.Range("E16").Value = "  -4.46%  "
$ws.Range("E17").Value = "  -2.12%  "
$ws.Range("E18").Value = "  +0.04%  "
$ws.Range("E19").Value = "  -1.09%  "
$ws.Range("E20").Value = "  +6.45%  "
$ws.Range("E21").Value = "  -0.77%  "
$ws.Range("E22").Value = "  -0.48%  "
$ws.Range("E23").Value = "  +0.60%  "
$ws.Range("E24").Value = "  -0.09%  "
$ws.Range("E25").Value = "  -0.96%  "
$ws.Range("E26").Value = "  +1.02%  "
$ws.Range("E27").Value = "  -0.51%  "
$ws.Range("E28").Value = "  -1.45%  "
$ws.Range("E29").Value = "  -1.10%  "
$ws.Range("E30").Value = "  -1.74%  "
$ws.Range("E31").Value = "  +0.68%  "
$ws.Range("E32").Value = "  -1.06%  "
$ws.Range("E33").Value = "  +0.02%  "
$ws.Range("E34").Value = "  -0.88%  "
$ws.Range("E35").Value = "  -1.11%  "
$ws.Range("E36").Value = "  -0.60%  "
$ws.Range("E37").Value = "  -1.87%  "
$ws.Range("E38").Value = "  +2.62%  "
$ws.Range("E39").Value = "  -0.86%  "
$ws.Range("E40").Value = "  -1.15%  "
$ws.Range("E41").Value = "  -1.97%  "
$ws.Range("E42").Value = "  +0.09%  "
$ws.Range("E43").Value = "  -0.17%  "
$ws.Range("E44").Value = "  -1.15%  "
$ws.Range("E45").Value = "  -0.98%  "
$ws.Range("E46").Value = "  -1.42%  "
$ws.Range("E47").Value = "  -2.98%  "
$ws.Range("E48").Value = "  -0.21%  "
$ws.Range("E49").Value = "  +0.28%  "
$ws.Range("E50").Value = "  -1.85%  "
$ws.Range("E51").Value = "  -0.42%  "

# Restore default (General) styling on column D now that values are entered
# as text, so the cells keep their original (no explicit style) appearance.
$priceRange.Style = "Normal"
